$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: fill in the new timesheet entry ---

# A5: date (reuse date style from A4 via copy/paste of formats only)
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A5").Value = 42722

# B5:D5: time values (reuse time style from B4:D4 via copy/paste of formats only)
$ws.Range("B4:D4").Copy() | Out-Null
$ws.Range("B5:D5").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$ws.Range("B5").Value = 0.4375
$ws.Range("C5").Value = 0.46388888888888885
$ws.Range("D5").Value = 0

# E5: comment text (reuse style from E4 via copy/paste of formats only)
$ws.Range("E4").Copy() | Out-Null
$ws.Range("E5").PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$ws.Range("E5").Value = "Implemented Jump, started working on dodge"

$excel.CutCopyMode = 0

# --- Update the selected/active cell shown in the worksheet view ---
$ws.Range("C6").Select() | Out-Null
